$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the three newly-added "Resumen de Reunión" entries in column H
$ws.Range("H9").Value = "Resumen de Reunión 19"
$ws.Range("H10").Value = "Resumen de Reunión 20"
$ws.Range("H11").Value = "Resumen de Reunión 21"

# Match the author's final selection state
$ws.Range("I10").Select()
